# "Generate Report for Handback" - mark a.md / b.md as handed back (in sync
# with en-US) for every locale sheet, filling in the Latest Target File /
# Latest Handback File / Latest Handback DateTime columns, and refresh the
# Status text wherever it is shown (Overview + each locale sheet).
#
# NOTE: existing cells (A/C/D columns, row 4, etc.) are intentionally left
# untouched so their original formatting/hyperlinks survive byte-for-byte;
# we only set values on the cells that actually change and add hyperlinks
# for the two brand-new linked cells per row (Latest Target File / Latest
# Handback File).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just refresh the Status text for both tracked files.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-locale sheets: zh-cn / de-de.
# Columns: A Source File Name, B Status, C Latest Handoff File,
#          D Latest Handoff Datetime, E Latest Target File,
#          F Latest Handback File, G Latest Handback DateTime,
#          H Handoff Reason, I Dependency From
# ---------------------------------------------------------------------

$locales = @(
    @{
        Sheet = "zh-cn"
        XlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9aa5d4b53d50a23e24a3803b001528a1a7ae6f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandbackTime = "2016-03-09 04:40:33"
    },
    @{
        Sheet = "de-de"
        XlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d16604b73abbd8b6dc8b497ec6e081e3c8715475/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandbackTime = "2016-03-09 04:40:38"
    }
)

$aUrl = "https://github.com/OpenLocalizationTest/oltest/blob/00a189f0462560d41d0ec36854be4f4f436927c5/e2e/a.md"

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # --- values: row 2 (a.md) and row 3 (b.md) get handed back ---
    $ws.Range("B2").Value = $newStatus
    $ws.Range("E2").Value = "a.md"
    $ws.Range("F2").Value = $locale.XlfName
    $ws.Range("G2").Value = $locale.HandbackTime

    $ws.Range("B3").Value = $newStatus
    $ws.Range("E3").Value = "a.md"
    $ws.Range("F3").Value = $locale.XlfName
    $ws.Range("G3").Value = $locale.HandbackTime

    # --- hyperlinks for the two brand-new linked cells per row ---
    $ws.Hyperlinks.Add($ws.Range("E2"), $aUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("F2"), $locale.XlfUrl, "", "", $locale.XlfName)
    $ws.Hyperlinks.Add($ws.Range("E3"), $aUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("F3"), $locale.XlfUrl, "", "", $locale.XlfName)
}
